# Update vendor quantities (previously commented out while testing the rest)
$wb = $excel.ActiveWorkbook

$grocery = $wb.Worksheets.Item("Grocery")
$pet     = $wb.Worksheets.Item("Pet")
$bath    = $wb.Worksheets.Item("Bath")

# --- Grocery sheet: update quantities in column B ---
$grocery.Range("B5").Value = 5
$grocery.Range("B6").Value = 5
$grocery.Range("B8").Value = 5

# --- Bath sheet: update quantities in column B ---
$bath.Range("B4").Value = 5
$bath.Range("B5").Value = 5
$bath.Range("B6").Value = 5

# --- Update selections on each sheet ---
$grocery.Range("B7").Select()
$pet.Range("B5").Select()
$bath.Range("B3").Select()

# --- Bath becomes the active (selected) sheet/tab ---
$bath.Activate()

$wb.Save()
